$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1012.63635
$ws.Range("I101").Value = 981
$ws.Range("J101").Value = 1039
$ws.Range("K101").Value = 2943
$ws.Range("L101").Value = 3117
$ws.Range("M101").Value = -1321
$ws.Range("N101").Value = -6361
$ws.Range("H112").Value = 3922.077
$ws.Range("J112").Value = 4453.364
$ws.Range("L112").Value = 13360.092
$ws.Range("N112").Value = -15576.092
$ws.Range("H131").Value = 3252229.8
$ws.Range("I131").Value = 6389.3335
$ws.Range("K131").Value = 19168.0005
$ws.Range("M131").Value = -14128.0005
$ws.Range("H137").Value = 10639886
$ws.Range("I137").Value = 16130122
$ws.Range("K137").Value = 48390366
$ws.Range("M137").Value = -48387816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8201.48
$ws.Range("I32").Value = 6771.8804
$ws.Range("K32").Value = 6771.8804
$ws.Range("M32").Value = -6484.8804
$ws.Range("H45").Value = 1545.7142
$ws.Range("I45").Value = 793.2727
$ws.Range("K45").Value = 793.2727
$ws.Range("M45").Value = -416.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3032646.2
$ws.Range("I134").Value = 2300.3872
$ws.Range("J134").Value = 50003010
$ws.Range("K134").Value = 6901.1616
$ws.Range("L134").Value = 150009030
$ws.Range("M134").Value = -4366.1616
$ws.Range("N134").Value = -150014100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3256.3076
$ws.Range("I58").Value = 2433.4
$ws.Range("K58").Value = 2433.4
$ws.Range("M58").Value = -2230.4
$ws.Range("H99").Value = 11636.68
$ws.Range("I99").Value = 5572
$ws.Range("K99").Value = 5572
$ws.Range("M99").Value = -4074
$ws.Range("H105").Value = 7352.3335
$ws.Range("I105").Value = 1174.1428
$ws.Range("J105").Value = 16001.8
$ws.Range("K105").Value = 1174.1428
$ws.Range("L105").Value = 16001.8
$ws.Range("M105").Value = 572.8571999999999
$ws.Range("N105").Value = -19495.8
$ws.Range("H126").Value = 11636.68
$ws.Range("I126").Value = 5572
$ws.Range("K126").Value = 16716
$ws.Range("M126").Value = -14246
$ws.Range("H132").Value = 2521.2632
$ws.Range("I132").Value = 2244.0625
$ws.Range("K132").Value = 6732.1875
$ws.Range("M132").Value = -4202.1875
$ws.Range("H134").Value = 1678.5264
$ws.Range("I134").Value = 1377.8572
$ws.Range("J134").Value = 2520.4
$ws.Range("K134").Value = 4133.571599999999
$ws.Range("L134").Value = 7561.200000000001
$ws.Range("M134").Value = -1598.571599999999
$ws.Range("N134").Value = -12631.2
$ws.Range("H136").Value = 3256.3076
$ws.Range("I136").Value = 2433.4
$ws.Range("K136").Value = 7300.200000000001
$ws.Range("M136").Value = -4750.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4554917
$ws.Range("J107").Value = 5692697.5
$ws.Range("L107").Value = 17078092.5
$ws.Range("N107").Value = -17081932.5
$ws.Range("H122").Value = 23898.357
$ws.Range("I122").Value = 33331.9
$ws.Range("J122").Value = 314.5
$ws.Range("K122").Value = 299987.1
$ws.Range("L122").Value = 2830.5
$ws.Range("M122").Value = -297537.1
$ws.Range("N122").Value = -7730.5
$ws.Range("H132").Value = 1748.6428
$ws.Range("I132").Value = 1919
$ws.Range("J132").Value = 1442
$ws.Range("K132").Value = 17271
$ws.Range("L132").Value = 12978
$ws.Range("M132").Value = -14741
$ws.Range("N132").Value = -18038
$ws.Range("H140").Value = 18753424
$ws.Range("I140").Value = 21429562
$ws.Range("K140").Value = 64288686
$ws.Range("M140").Value = -64283506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12265.429
$ws.Range("I70").Value = 11590.444
$ws.Range("J70").Value = 12417.3
$ws.Range("K70").Value = 11590.444
$ws.Range("L70").Value = 12417.3
$ws.Range("M70").Value = -11320.444
$ws.Range("N70").Value = -12957.3
$ws.Range("H73").Value = 12265.429
$ws.Range("I73").Value = 11590.444
$ws.Range("J73").Value = 12417.3
$ws.Range("K73").Value = 11590.444
$ws.Range("L73").Value = 12417.3
$ws.Range("M73").Value = -10654.444
$ws.Range("N73").Value = -14289.3
$ws.Range("H80").Value = 200002990
$ws.Range("J80").Value = 200002990
$ws.Range("L80").Value = 200002990
$ws.Range("N80").Value = -200004986
$ws.Range("H83").Value = 200002990
$ws.Range("J83").Value = 200002990
$ws.Range("L83").Value = 1000014950
$ws.Range("N83").Value = -1000024934
$ws.Range("H122").Value = 5431.3228
$ws.Range("I122").Value = 3946.2173
$ws.Range("K122").Value = 11838.6519
$ws.Range("M122").Value = -9388.651899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7835.1763
$ws.Range("I7").Value = 7626.25
$ws.Range("J7").Value = 8020.8887
$ws.Range("K7").Value = 7626.25
$ws.Range("L7").Value = 8020.8887
$ws.Range("M7").Value = -7514.25
$ws.Range("N7").Value = -8244.8887
$ws.Range("H46").Value = 816.5
$ws.Range("I46").Value = 495.81818
$ws.Range("K46").Value = 495.81818
$ws.Range("M46").Value = -307.81818
$ws.Range("H55").Value = 822.2759
$ws.Range("J55").Value = 1291.1333
$ws.Range("L55").Value = 1291.1333
$ws.Range("N55").Value = -1637.1333
$ws.Range("H61").Value = 2537.1428
$ws.Range("I61").Value = 2229.2856
$ws.Range("J61").Value = 3768.5715
$ws.Range("K61").Value = 2229.2856
$ws.Range("L61").Value = 3768.5715
$ws.Range("M61").Value = -2027.2856
$ws.Range("N61").Value = -4172.5715
$ws.Range("H82").Value = 4059.5454
$ws.Range("I82").Value = 1457.125
$ws.Range("K82").Value = 1457.125
$ws.Range("M82").Value = -1096.125
$ws.Range("H85").Value = 4059.5454
$ws.Range("I85").Value = 1457.125
$ws.Range("K85").Value = 1457.125
$ws.Range("M85").Value = -209.125
$ws.Range("H113").Value = 2537.1428
$ws.Range("I113").Value = 2229.2856
$ws.Range("J113").Value = 3768.5715
$ws.Range("K113").Value = 2229.2856
$ws.Range("L113").Value = 3768.5715
$ws.Range("M113").Value = -59.28560000000016
$ws.Range("N113").Value = -8108.5715
$ws.Range("H122").Value = 3994.641
$ws.Range("I122").Value = 3494.3057
$ws.Range("K122").Value = 10482.9171
$ws.Range("M122").Value = -8032.917099999999
$ws.Range("H126").Value = 7835.1763
$ws.Range("I126").Value = 7626.25
$ws.Range("J126").Value = 8020.8887
$ws.Range("K126").Value = 22878.75
$ws.Range("L126").Value = 24062.6661
$ws.Range("M126").Value = -20408.75
$ws.Range("N126").Value = -29002.6661
$ws.Range("H132").Value = 3840.75
$ws.Range("I132").Value = 2021.56
$ws.Range("K132").Value = 6064.68
$ws.Range("M132").Value = -3534.68
$ws.Range("H136").Value = 3702.2307
$ws.Range("I136").Value = 1876.1111
$ws.Range("J136").Value = 7811
$ws.Range("K136").Value = 5628.3333
$ws.Range("L136").Value = 23433
$ws.Range("M136").Value = -3078.3333
$ws.Range("N136").Value = -28533

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 15180
$ws.Range("H113").Value = 400.66666
$ws.Range("I113").Value = 215.25
$ws.Range("J113").Value = 771.5
$ws.Range("K113").Value = 645.75
$ws.Range("L113").Value = 2314.5
$ws.Range("M113").Value = 1524.25
$ws.Range("N113").Value = -6654.5
$ws.Range("H136").Value = 288598.84
$ws.Range("I136").Value = 2685.3667
$ws.Range("K136").Value = 8056.1001
$ws.Range("M136").Value = -5506.1001
